$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "Lösungsvorschlag" value in B5
$ws.Range("B5").Value = "Bootvorgang auf Netzwerk und HDD beschränken"

# Add a new row with the justification ("Begründung")
$ws.Range("A6").Value = "Begründung"
$ws.Range("B6").Value = "Zum Schutz vor Fremdbenutzung"

# Keep selection consistent with the saved file (active cell moves to B6)
$ws.Range("B6").Select()
